$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 638
$ws1.Range("F8").Value = 1401
$ws1.Range("F9").Value = 4020
$ws1.Range("F10").Value = 88

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 57

# Sheet "全部类型" (all types, aggregated)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 638
$ws4.Range("F8").Value = 1401
$ws4.Range("F9").Value = 4020
$ws4.Range("F10").Value = 88
$ws4.Range("F11").Value = 57
